# Auto-generated Excel COM-interop script
# Applies the weekly CompStat data refresh described by the commit:
#   "New crime data collected"
#
# - Updates the report title/volume and the covered week dates (shared strings)
# - Updates ~195 numeric statistics cells (counts + computed % changes) in rows 14-46

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number and reporting week dates ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Data table refresh (rows 14-46) ---

# Row 14
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N14").Value = -100
$ws.Range("N14").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = -15.384615384615
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = -25
$ws.Range("M15").Value = -25
$ws.Range("M15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = -57.142857142857
# Row 16
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 27
$ws.Range("E16").Value = -7.407407407407
$ws.Range("F16").Value = 94
$ws.Range("G16").Value = 120
$ws.Range("H16").Value = -21.666666666666
$ws.Range("I16").Value = 33
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = -21.428571428571
$ws.Range("L16").Value = -34
$ws.Range("M16").Value = -26.666666666666
$ws.Range("N16").Value = -89.655172413793
# Row 17
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 91.666666666666
$ws.Range("F17").Value = 178
$ws.Range("G17").Value = 160
$ws.Range("H17").Value = 11.25
$ws.Range("I17").Value = 83
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 50.909090909090
$ws.Range("L17").Value = 33.870967741935
$ws.Range("M17").Value = 124.324324324324
$ws.Range("N17").Value = -27.192982456140
# Row 18
$ws.Range("C18").Value = 35
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = -14.634146341463
$ws.Range("G18").Value = 174
$ws.Range("H18").Value = -21.264367816092
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = -20.967741935483
$ws.Range("M18").Value = -39.506172839506
$ws.Range("N18").Value = -87.626262626262
# Row 19
$ws.Range("C19").Value = 185
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = -7.5
$ws.Range("F19").Value = 740
$ws.Range("G19").Value = 818
$ws.Range("H19").Value = -9.535452322738
$ws.Range("I19").Value = 310
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = -7.185628742514
$ws.Range("L19").Value = -15.300546448087
$ws.Range("M19").Value = 4.026845637583
$ws.Range("N19").Value = -69.061876247505
# Row 20
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 6.666666666666
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 6
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -75
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -98.067632850241
# Row 21
$ws.Range("C21").Value = 296
$ws.Range("D21").Value = 297
$ws.Range("E21").Value = -0.336700336700
$ws.Range("F21").Value = 1177
$ws.Range("G21").Value = 1301
$ws.Range("H21").Value = -9.531129900076
$ws.Range("I21").Value = 482
$ws.Range("J21").Value = 505
$ws.Range("K21").Value = -4.554455445544
$ws.Range("L21").Value = -14.081996434937
$ws.Range("M21").Value = 2.991452991452
$ws.Range("N21").Value = -76.46484375
# Row 22
$ws.Range("C22").Value = 11
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 48
$ws.Range("G22").Value = 52
$ws.Range("H22").Value = -7.692307692307
$ws.Range("I22").Value = 19
$ws.Range("J22").Value = 14
$ws.Range("K22").Value = 35.714285714285
$ws.Range("L22").Value = -13.636363636363
$ws.Range("M22").Value = -9.523809523809
# Row 23
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 133.333333333333
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 37.5
$ws.Range("L23").Value = 22.222222222222
$ws.Range("M23").Value = 10
# Row 24
$ws.Range("C24").Value = 380
$ws.Range("D24").Value = 373
$ws.Range("E24").Value = 1.876675603217
$ws.Range("F24").Value = 1373
$ws.Range("G24").Value = 1507
$ws.Range("H24").Value = -8.891838088918
$ws.Range("I24").Value = 498
$ws.Range("J24").Value = 552
$ws.Range("K24").Value = -9.782608695652
$ws.Range("L24").Value = -16.442953020134
$ws.Range("M24").Value = 16.083916083916
# Row 25
$ws.Range("C25").Value = 291
$ws.Range("D25").Value = 290
$ws.Range("E25").Value = 0.344827586206
$ws.Range("F25").Value = 1062
$ws.Range("G25").Value = 1233
$ws.Range("H25").Value = -13.868613138686
$ws.Range("I25").Value = 399
$ws.Range("J25").Value = 431
$ws.Range("K25").Value = -7.424593967517
$ws.Range("L25").Value = -21.917808219178
# Row 26
$ws.Range("C26").Value = 86
$ws.Range("D26").Value = 81
$ws.Range("E26").Value = 6.172839506172
$ws.Range("F26").Value = 353
$ws.Range("G26").Value = 343
$ws.Range("H26").Value = 2.915451895043
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 126
$ws.Range("K26").Value = 16.666666666666
$ws.Range("L26").Value = 8.888888888888
$ws.Range("M26").Value = 54.736842105263
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = -21.428571428571
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -57.142857142857
# Row 28
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = 71.428571428571
$ws.Range("F28").Value = 51
$ws.Range("G28").Value = 49
$ws.Range("H28").Value = 4.081632653061
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 145.454545454545
$ws.Range("L28").Value = 28.571428571428
# Row 29
$ws.Range("G29").Value = 1
$ws.Range("L29").Value = -100
$ws.Range("L29").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 30
$ws.Range("G30").Value = 1
$ws.Range("L30").Value = -100
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 31
$ws.Range("D31").Value = 2
$ws.Range("F31").Value = 5
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 3
$ws.Range("L31").Value = -100
$ws.Range("L31").NumberFormat = '#,##0.0;"-"#,##0.0'
# Row 33
$ws.Range("F33").Value = 2
# Row 41
$ws.Range("J41").Value = 1615
$ws.Range("K41").Value = -40.207330618289
$ws.Range("L41").Value = -61.483424755545
$ws.Range("M41").Value = -84.726688102893
$ws.Range("N41").Value = -89.136284138302
# Row 42
$ws.Range("J42").Value = 2420
$ws.Range("K42").Value = 41.190198366394
$ws.Range("L42").Value = 6.093818500657
$ws.Range("M42").Value = -30.399769916594
$ws.Range("N42").Value = -39.454590943207
# Row 43
$ws.Range("J43").Value = 1981
$ws.Range("K43").Value = -46.747311827957
$ws.Range("L43").Value = -62.805107022155
$ws.Range("M43").Value = -83.772935779816
$ws.Range("N43").Value = -87.688004972032
# Row 44
$ws.Range("J44").Value = 10254
$ws.Range("K44").Value = -38.499370239309
$ws.Range("L44").Value = -49.063633202523
$ws.Range("M44").Value = -69.042659179422
$ws.Range("N44").Value = -77.117225681194
# Row 46
$ws.Range("J46").Value = 16836
$ws.Range("K46").Value = -36.183761655674
$ws.Range("L46").Value = -51.336821111657
$ws.Range("M46").Value = -74.469246633507
$ws.Range("N46").Value = -81.197440278755
